$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 1229.52
$ws.Cells.Item(6, 9).Value = 447.33334
$ws.Cells.Item(6, 10).Value = 20002
$ws.Cells.Item(6, 11).Value = 1342.00002
$ws.Cells.Item(6, 12).Value = 60006
$ws.Cells.Item(6, 13).Value = -1230.00002
$ws.Cells.Item(6, 14).Value = -60230
$ws.Cells.Item(38, 8).Value = 625377.7
$ws.Cells.Item(38, 9).Value = 909181.6
$ws.Cells.Item(38, 10).Value = 1009
$ws.Cells.Item(38, 11).Value = 2727544.8
$ws.Cells.Item(38, 12).Value = 3027
$ws.Cells.Item(38, 13).Value = -2727172.8
$ws.Cells.Item(38, 14).Value = -3771
$ws.Cells.Item(137, 8).Value = 22539.256
$ws.Cells.Item(137, 9).Value = 24647.785
$ws.Cells.Item(137, 11).Value = 73943.355
$ws.Cells.Item(137, 13).Value = -71393.355
$ws.Cells.Item(141, 8).Value = 2646.2666
$ws.Cells.Item(141, 9).Value = 1854.8889
$ws.Cells.Item(141, 10).Value = 3833.3333
$ws.Cells.Item(141, 11).Value = 5564.6667
$ws.Cells.Item(141, 12).Value = 11499.9999
$ws.Cells.Item(141, 13).Value = -384.6666999999998
$ws.Cells.Item(141, 14).Value = -21859.9999
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1516.0358
$ws.Cells.Item(61, 9).Value = 793.1818
$ws.Cells.Item(61, 10).Value = 4166.5
$ws.Cells.Item(61, 11).Value = 793.1818
$ws.Cells.Item(61, 12).Value = 4166.5
$ws.Cells.Item(61, 13).Value = -581.1818
$ws.Cells.Item(61, 14).Value = -4590.5
$ws.Cells.Item(74, 8).Value = 37747.637
$ws.Cells.Item(74, 9).Value = 61976.547
$ws.Cells.Item(74, 10).Value = 1404.2727
$ws.Cells.Item(74, 11).Value = 61976.547
$ws.Cells.Item(74, 12).Value = 1404.2727
$ws.Cells.Item(74, 13).Value = -61102.547
$ws.Cells.Item(74, 14).Value = -3152.2727
$ws.Cells.Item(77, 8).Value = 37747.637
$ws.Cells.Item(77, 9).Value = 61976.547
$ws.Cells.Item(77, 10).Value = 1404.2727
$ws.Cells.Item(77, 11).Value = 309882.735
$ws.Cells.Item(77, 12).Value = 7021.363499999999
$ws.Cells.Item(77, 13).Value = -305514.735
$ws.Cells.Item(77, 14).Value = -15757.3635
$ws.Cells.Item(101, 8).Value = 39998
$ws.Cells.Item(101, 10).Value = 39998
$ws.Cells.Item(101, 12).Value = 39998
$ws.Cells.Item(101, 14).Value = -46488
$ws.Cells.Item(132, 8).Value = 2438593.2
$ws.Cells.Item(132, 9).Value = 2915651.2
$ws.Cells.Item(132, 10).Value = 920682.0600000001
$ws.Cells.Item(132, 11).Value = 8746953.600000001
$ws.Cells.Item(132, 12).Value = 2762046.18
$ws.Cells.Item(132, 13).Value = -8744423.600000001
$ws.Cells.Item(132, 14).Value = -2767106.18
$ws.Cells.Item(136, 8).Value = 1516.0358
$ws.Cells.Item(136, 9).Value = 793.1818
$ws.Cells.Item(136, 10).Value = 4166.5
$ws.Cells.Item(136, 11).Value = 2379.5454
$ws.Cells.Item(136, 12).Value = 12499.5
$ws.Cells.Item(136, 13).Value = 170.4546
$ws.Cells.Item(136, 14).Value = -17599.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1647.5
$ws.Cells.Item(99, 9).Value = 1196.6666
$ws.Cells.Item(99, 10).Value = 3000
$ws.Cells.Item(99, 11).Value = 1196.6666
$ws.Cells.Item(99, 12).Value = 3000
$ws.Cells.Item(99, 13).Value = 301.3334
$ws.Cells.Item(99, 14).Value = -5996
$ws.Cells.Item(134, 8).Value = 56353.473
$ws.Cells.Item(134, 9).Value = 2423.6365
$ws.Cells.Item(134, 10).Value = 130507
$ws.Cells.Item(134, 11).Value = 7270.9095
$ws.Cells.Item(134, 12).Value = 391521
$ws.Cells.Item(134, 13).Value = -4735.9095
$ws.Cells.Item(134, 14).Value = -396591
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 31945.727
$ws.Cells.Item(31, 9).Value = 53068.695
$ws.Cells.Item(31, 10).Value = 8811.048000000001
$ws.Cells.Item(31, 11).Value = 53068.695
$ws.Cells.Item(31, 12).Value = 8811.048000000001
$ws.Cells.Item(31, 13).Value = -52773.695
$ws.Cells.Item(31, 14).Value = -9401.048000000001
$ws.Cells.Item(34, 8).Value = 31945.727
$ws.Cells.Item(34, 9).Value = 53068.695
$ws.Cells.Item(34, 10).Value = 8811.048000000001
$ws.Cells.Item(34, 11).Value = 53068.695
$ws.Cells.Item(34, 12).Value = 8811.048000000001
$ws.Cells.Item(34, 13).Value = -52866.695
$ws.Cells.Item(34, 14).Value = -9215.048000000001
$ws.Cells.Item(35, 8).Value = 30658.428
$ws.Cells.Item(35, 9).Value = 1131.25
$ws.Cells.Item(35, 10).Value = 70028
$ws.Cells.Item(35, 11).Value = 1131.25
$ws.Cells.Item(35, 12).Value = 70028
$ws.Cells.Item(35, 13).Value = -837.25
$ws.Cells.Item(35, 14).Value = -70616
$ws.Cells.Item(51, 8).Value = 8499.885
$ws.Cells.Item(51, 10).Value = 8499.885
$ws.Cells.Item(51, 12).Value = 8499.885
$ws.Cells.Item(51, 14).Value = -9971.885
$ws.Cells.Item(58, 8).Value = 8080.5713
$ws.Cells.Item(58, 9).Value = 1344.4445
$ws.Cells.Item(58, 10).Value = 20205.6
$ws.Cells.Item(58, 11).Value = 1344.4445
$ws.Cells.Item(58, 12).Value = 20205.6
$ws.Cells.Item(58, 13).Value = -1141.4445
$ws.Cells.Item(58, 14).Value = -20611.6
$ws.Cells.Item(61, 8).Value = 8499.885
$ws.Cells.Item(61, 10).Value = 8499.885
$ws.Cells.Item(61, 12).Value = 8499.885
$ws.Cells.Item(61, 14).Value = -9195.885
$ws.Cells.Item(106, 8).Value = 58742.5
$ws.Cells.Item(106, 10).Value = 58742.5
$ws.Cells.Item(106, 12).Value = 58742.5
$ws.Cells.Item(106, 14).Value = -61266.5
$ws.Cells.Item(132, 8).Value = 2361.2163
$ws.Cells.Item(132, 9).Value = 1872.125
$ws.Cells.Item(132, 10).Value = 3264.1538
$ws.Cells.Item(132, 11).Value = 5616.375
$ws.Cells.Item(132, 12).Value = 9792.4614
$ws.Cells.Item(132, 13).Value = -3086.375
$ws.Cells.Item(132, 14).Value = -14852.4614
$ws.Cells.Item(134, 8).Value = 14707582
$ws.Cells.Item(134, 9).Value = 1487.579
$ws.Cells.Item(134, 10).Value = 33335300
$ws.Cells.Item(134, 11).Value = 4462.737
$ws.Cells.Item(134, 12).Value = 100005900
$ws.Cells.Item(134, 13).Value = -1927.737
$ws.Cells.Item(134, 14).Value = -100010970
$ws.Cells.Item(136, 8).Value = 8080.5713
$ws.Cells.Item(136, 9).Value = 1344.4445
$ws.Cells.Item(136, 10).Value = 20205.6
$ws.Cells.Item(136, 11).Value = 4033.3335
$ws.Cells.Item(136, 12).Value = 60616.8
$ws.Cells.Item(136, 13).Value = -1483.3335
$ws.Cells.Item(136, 14).Value = -65716.79999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 3038.1
$ws.Cells.Item(3, 9).Value = 1485.8334
$ws.Cells.Item(3, 10).Value = 4072.9443
$ws.Cells.Item(3, 11).Value = 4457.5002
$ws.Cells.Item(3, 12).Value = 12218.8329
$ws.Cells.Item(3, 13).Value = -4345.5002
$ws.Cells.Item(3, 14).Value = -12442.8329
$ws.Cells.Item(49, 8).Value = 2000.4445
$ws.Cells.Item(49, 10).Value = 2000.4445
$ws.Cells.Item(49, 12).Value = 6001.333500000001
$ws.Cells.Item(49, 14).Value = -6313.333500000001
$ws.Cells.Item(113, 8).Value = 622.7619
$ws.Cells.Item(113, 9).Value = 448
$ws.Cells.Item(113, 11).Value = 1344
$ws.Cells.Item(113, 13).Value = 826
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 85325.836
$ws.Cells.Item(132, 9).Value = 2172.3635
$ws.Cells.Item(132, 10).Value = 1000014
$ws.Cells.Item(132, 11).Value = 6517.0905
$ws.Cells.Item(132, 12).Value = 3000042
$ws.Cells.Item(132, 13).Value = -3987.0905
$ws.Cells.Item(132, 14).Value = -3005102
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(35, 8).Value = 14301
$ws.Cells.Item(35, 9).Value = 2168.3333
$ws.Cells.Item(35, 10).Value = 32500
$ws.Cells.Item(35, 11).Value = 2168.3333
$ws.Cells.Item(35, 12).Value = 32500
$ws.Cells.Item(35, 13).Value = -1832.3333
$ws.Cells.Item(35, 14).Value = -33172
$ws.Cells.Item(104, 8).Value = 38322
$ws.Cells.Item(104, 10).Value = 38322
$ws.Cells.Item(104, 12).Value = 38322
$ws.Cells.Item(104, 14).Value = -45310
$ws.Cells.Item(132, 8).Value = 302748.56
$ws.Cells.Item(132, 9).Value = 76335.336
$ws.Cells.Item(132, 10).Value = 772991.4
$ws.Cells.Item(132, 11).Value = 229006.008
$ws.Cells.Item(132, 12).Value = 2318974.2
$ws.Cells.Item(132, 13).Value = -226476.008
$ws.Cells.Item(132, 14).Value = -2324034.2
$ws.Cells.Item(136, 8).Value = 418168.88
$ws.Cells.Item(136, 9).Value = 667466.9
$ws.Cells.Item(136, 10).Value = 2672.2222
$ws.Cells.Item(136, 11).Value = 2002400.7
$ws.Cells.Item(136, 12).Value = 8016.6666
$ws.Cells.Item(136, 13).Value = -1999850.7
$ws.Cells.Item(136, 14).Value = -13116.6666
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(101, 8).Value = 14254.182
$ws.Cells.Item(101, 10).Value = 14254.182
$ws.Cells.Item(101, 12).Value = 14254.182
$ws.Cells.Item(101, 14).Value = -20744.182
$ws.Cells.Item(122, 8).Value = 4013.7454
$ws.Cells.Item(122, 9).Value = 3666.3
$ws.Cells.Item(122, 10).Value = 4940.2666
$ws.Cells.Item(122, 11).Value = 10998.9
$ws.Cells.Item(122, 12).Value = 14820.7998
$ws.Cells.Item(122, 13).Value = -8548.900000000001
$ws.Cells.Item(122, 14).Value = -19720.7998
$ws.Cells.Item(132, 8).Value = 4355.727
$ws.Cells.Item(132, 9).Value = 1111.3793
$ws.Cells.Item(132, 11).Value = 3334.1379
$ws.Cells.Item(132, 13).Value = -804.1379000000002
$ws.Cells.Item(136, 8).Value = 4287963.5
$ws.Cells.Item(136, 9).Value = 4204019
$ws.Cells.Item(136, 10).Value = 5001492.5
$ws.Cells.Item(136, 11).Value = 12612057
$ws.Cells.Item(136, 12).Value = 15004477.5
$ws.Cells.Item(136, 13).Value = -12609507
$ws.Cells.Item(136, 14).Value = -15009577.5
